$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (4th repetition of the scale pattern)
$data = @(
    @("c3_4", 130, 2),
    @("d3_4", 146, 2),
    @("e3_4", 164, 2),
    @("f3_4", 174, 2),
    @("g3_4", 196, 2),
    @("a3_4", 208, 2),
    @("b3_4", 246, 2),
    @("c4_4", 261, 2)
)

$row = 26
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

# Update the view: select A34 (the next empty row below the new data) and
# scroll the window so row 15 is at the top, matching Excel's post-entry view.
$ws.Range("A34").Select()
$excel.ActiveWindow.ScrollRow = 15
